# Generate Report for Handback
#
# The 21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.md file has been handed back
# (translations for zh-cn and de-de are now in sync with en-US). Update the
# localization status report: flip its status from "Ready for handoff" to
# "Handed back: in sync with en-US" on the Overview sheet and on each
# language sheet, and record the handback details (target file, handback
# file + link, handback datetime) on the language sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: row 2 is the 21ed720f file, columns B (zh-cn) / C (de-de) ---
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack

# Colour used by this workbook's existing hyperlink cells (RGB 6495ED).
$hyperlinkColor = 15570276

function Set-HandbackRow($ws, $mdAddress, $mdDisplay, $xlfAddress, $xlfDisplay, $handbackDateTime) {
    # Status -> handed back, in sync with en-US
    $ws.Range("C2").Value = $statusHandedBack

    # Latest Target File (F2) - the handed-back source file
    $ws.Hyperlinks.Add($ws.Range("F2"), $mdAddress, [Type]::Missing, [Type]::Missing, $mdDisplay) | Out-Null
    $ws.Range("F2").Font.Underline = 2
    $ws.Range("F2").Font.Color = $hyperlinkColor

    # Latest Handback File (G2) - the handed-back translation file
    $ws.Hyperlinks.Add($ws.Range("G2"), $xlfAddress, [Type]::Missing, [Type]::Missing, $xlfDisplay) | Out-Null
    $ws.Range("G2").Font.Underline = 2
    $ws.Range("G2").Font.Color = $hyperlinkColor

    # Latest Handback DateTime (H2)
    $ws.Range("H2").Value = $handbackDateTime
}

# --- zh-cn sheet ---
Set-HandbackRow $wsZhCn `
    "https://github.com/OpenLocalizationTest/oltest/blob/46f9ef041379da3de779195861d7389e2ee84d21/e2e/21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.md" `
    "21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce26147698f3f0b826483f80b87af99114e7ad8e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.7b32db215b5030fc8eef5443a292995986e0f93c.zh-cn.xlf" `
    "21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.7b32db215b5030fc8eef5443a292995986e0f93c.zh-cn.xlf" `
    "2016-03-11 08:22:13"

# --- de-de sheet ---
Set-HandbackRow $wsDeDe `
    "https://github.com/OpenLocalizationTest/oltest/blob/46f9ef041379da3de779195861d7389e2ee84d21/e2e/21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.md" `
    "21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2fd70ef8883cb09e7e1f017b2162b9febc503673/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.7b32db215b5030fc8eef5443a292995986e0f93c.de-de.xlf" `
    "21ed720f-7e9e-4c8d-bada-a09f64d3d4ac.7b32db215b5030fc8eef5443a292995986e0f93c.de-de.xlf" `
    "2016-03-11 08:22:18"
